# Update countries & provincias Spain
# Refresh the COVID-19 country statistics on the "Pais" sheet with the
# latest reported figures (total cases, new cases, active cases,
# recovered, critical cases, new deaths, total deaths).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 7960416
$ws.Cells.Item(4, 3).Value = 14911
$ws.Cells.Item(4, 4).Value = 5093861
$ws.Cells.Item(4, 5).Value = 2647191
$ws.Cells.Item(4, 7).Value = 82
$ws.Cells.Item(4, 8).Value = 219364

# Row 5 - India
$ws.Cells.Item(5, 2).Value = 7073958
$ws.Cells.Item(5, 3).Value = 22415
$ws.Cells.Item(5, 4).Value = 6098259
$ws.Cells.Item(5, 5).Value = 867176
$ws.Cells.Item(5, 7).Value = 152
$ws.Cells.Item(5, 8).Value = 108523

# Row 6 - Brasil
$ws.Cells.Item(6, 2).Value = 5095586
$ws.Cells.Item(6, 3).Value = 3746
$ws.Cells.Item(6, 5).Value = 491526
$ws.Cells.Item(6, 7).Value = 102
$ws.Cells.Item(6, 8).Value = 150338

# Row 23 - Turquia
$ws.Cells.Item(23, 2).Value = 335533
$ws.Cells.Item(23, 3).Value = 1502
$ws.Cells.Item(23, 4).Value = 294357
$ws.Cells.Item(23, 5).Value = 32339
$ws.Cells.Item(23, 7).Value = 59
$ws.Cells.Item(23, 8).Value = 8837

# Row 25 - Alemania
$ws.Cells.Item(25, 2).Value = 324938
$ws.Cells.Item(25, 3).Value = 1485
$ws.Cells.Item(25, 5).Value = 41739
$ws.Cells.Item(25, 7).Value = 8
$ws.Cells.Item(25, 8).Value = 9699

# Row 29 - Canada
$ws.Cells.Item(29, 2).Value = 181772
$ws.Cells.Item(29, 3).Value = 1593
$ws.Cells.Item(29, 4).Value = 153206
$ws.Cells.Item(29, 5).Value = 18954
$ws.Cells.Item(29, 7).Value = 4
$ws.Cells.Item(29, 8).Value = 9612

# Row 40 - Chequia
$ws.Cells.Item(40, 2).Value = 114547
$ws.Cells.Item(40, 3).Value = 542
$ws.Cells.Item(40, 4).Value = 54980
$ws.Cells.Item(40, 5).Value = 58585
$ws.Cells.Item(40, 7).Value = 34
$ws.Cells.Item(40, 8).Value = 982

# Row 87 - Grecia
$ws.Cells.Item(87, 2).Value = 22358
$ws.Cells.Item(87, 3).Value = 280
$ws.Cells.Item(87, 5).Value = 11920
$ws.Cells.Item(87, 7).Value = 13
$ws.Cells.Item(87, 8).Value = 449

# Row 109 - Mozambique
$ws.Cells.Item(109, 2).Value = 10001
$ws.Cells.Item(109, 3).Value = 157
$ws.Cells.Item(109, 4).Value = 7338
$ws.Cells.Item(109, 5).Value = 2592
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 71

# Row 135 - Sri Lanka
$ws.Cells.Item(135, 2).Value = 4752
$ws.Cells.Item(135, 3).Value = 124
$ws.Cells.Item(135, 5).Value = 1432

# Row 156 - Sierra Leona
$ws.Cells.Item(156, 2).Value = 2306
$ws.Cells.Item(156, 3).Value = 6
$ws.Cells.Item(156, 4).Value = 1736
$ws.Cells.Item(156, 5).Value = 498

# Row 160 - Republica de Chipre
$ws.Cells.Item(160, 2).Value = 2006
$ws.Cells.Item(160, 3).Value = 20
$ws.Cells.Item(160, 5).Value = 537

# Row 174 - Curazao
$ws.Cells.Item(174, 2).Value = 583
$ws.Cells.Item(174, 3).Value = 12
$ws.Cells.Item(174, 4).Value = 315
$ws.Cells.Item(174, 5).Value = 267
